# Product Backlog update
# - Sprint # bump for row 6 (D6)
# - Fill in Assigned To / Start / Finish for several user stories
# - Update Status for several rows
# - Move selection to O14 (last user interaction before save)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Agile Product Backlog")

# Row 6: sprint number change
$ws.Range("D6").Value = 2

# Row 10
$ws.Range("F10").Value = "27th Jan"
$ws.Range("G10").Value = "30th Jan"
$ws.Range("K10").Value = "Completed"

# Row 13
$ws.Range("E13").Value = "Ross"
$ws.Range("F13").Value = "27th Jan"
$ws.Range("G13").Value = "29th Jan"
$ws.Range("K13").Value = "Completed"

# Row 14
$ws.Range("E14").Value = "Kamila"
$ws.Range("F14").Value = "27th Jan"
$ws.Range("G14").Value = "28th Jan"
$ws.Range("K14").Value = "Completed"

# Row 15
$ws.Range("F15").Value = "28th Jan"
$ws.Range("K15").Value = "In Progress"

# Row 24
$ws.Range("E24").Value = "Kayla"
$ws.Range("F24").Value = "27th Jan"
$ws.Range("G24").Value = "31st Jan"
$ws.Range("K24").Value = "Completed"

# Final cursor position left by the author before saving
$ws.Range("O14").Select()
